# Adapt templates to work with latest release
# - bump the template's own version marker
# - bump the "last tested with" ReportServer build reference
# - move the active selection down one row (A3 -> A4)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Customers")

$ws.Range("A4").Value = "(Version: 1.0.1)"
$ws.Range("A5").Value = "(Last tested with: ReportServer 4.0.0-6053) "

$ws.Activate()
$ws.Range("A4").Select()
